$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 9 (old articles), keeping only header + row 2
$ws.Range("A3:E9").EntireRow.Delete()

# Update row 2 with the new article data
$ws.Range("A2").Value = "Diam-Diam Prabowo Pantau Kinerja Menkeu Purbaya via Medsos"
$ws.Range("B2").Value = "2025-10-01T07:29:36+07:00"
$ws.Range("C2").Value = "Tidak Diketahui"
$ws.Range("D2").Value = "https://www.liputan6.com/news/read/6172970/diam-diam-prabowo-pantau-kinerja-menkeu-purbaya-via-medsos"
$ws.Range("E2").Value = "purbaya"
